$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "69.621.55"
Set-TextValue "E2" "  -0.84%  "
Set-TextValue "D3" "2.499.38"
Set-TextValue "E3" "  -0.89%  "
Set-TextValue "E4" "  +0.07%  "
Set-TextValue "D5" "572.63"
Set-TextValue "E5" "  -1.21%  "
Set-TextValue "D6" "166.07"
Set-TextValue "E6" "  -1.40%  "
Set-TextValue "E7" "  -0.01%  "
Set-TextValue "D8" "0.513"
Set-TextValue "E8" "  -1.09%  "
Set-TextValue "D9" "2.500.16"
Set-TextValue "E9" "  -0.81%  "
Set-TextValue "E10" "  -1.81%  "
Set-TextValue "D11" "0.167"
Set-TextValue "E11" "  -0.23%  "
Set-TextValue "D12" "0.357"
Set-TextValue "E12" "  +2.17%  "
Set-TextValue "E13" "  +0.66%  "
Set-TextValue "D14" "2.957.53"
Set-TextValue "E14" "  -0.91%  "
Set-TextValue "D15" "69.526.81"
Set-TextValue "E15" "  -0.82%  "
Set-TextValue "E16" "  -0.08%  "
Set-TextValue "D17" "24.64"
Set-TextValue "E17" "  -2.28%  "
Set-TextValue "D18" "2.499.77"
Set-TextValue "E18" "  -0.40%  "
Set-TextValue "D19" "11.16"
Set-TextValue "E19" "  -1.86%  "
Set-TextValue "D20" "7.41"
Set-TextValue "E20" "  -5.59%  "
Set-TextValue "D21" "347.90"
Set-TextValue "E21" "  -1.26%  "
Set-TextValue "D22" "3.89"
Set-TextValue "E22" "  -1.14%  "
Set-TextValue "E23" "  -1.53%  "
Set-TextValue "E24" "  -0.05%  "
Set-TextValue "D25" "70.72"
Set-TextValue "E25" "  +1.62%  "
Set-TextValue "D26" "3.90"
Set-TextValue "E26" "  -2.85%  "
Set-TextValue "B27" "WrappedeETH"
Set-TextValue "C27" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D27" "2.628.69"
Set-TextValue "E27" "  -1.47%  "
Set-TextValue "B28" "Aptos"
Set-TextValue "C28" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D28" "8.69"
Set-TextValue "E28" "  -3.95%  "
Set-TextValue "D29" "1.00"
Set-TextValue "E29" "  +1.22%  "
Set-TextValue "E30" "  -2.79%  "
Set-TextValue "D31" "7.83"
Set-TextValue "E31" "  -0.91%  "
Set-TextValue "D32" "455.02"
Set-TextValue "E32" "  -2.14%  "
Set-TextValue "E33" "  -6.26%  "
Set-TextValue "D34" "1.72"
Set-TextValue "E34" "  -1.80%  "
Set-TextValue "E35" "  +0.01%  "
Set-TextValue "B36" "Kaspa"
Set-TextValue "C36" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D36" "0.115"
Set-TextValue "E36" "  -3.39%  "
Set-TextValue "B37" "Monero"
Set-TextValue "C37" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D37" "156.74"
Set-TextValue "E37" "  +1.59%  "
Set-TextValue "E38" "  +0.13%  "
Set-TextValue "D39" "18.31"
Set-TextValue "E39" "  -0.99%  "
Set-TextValue "E40" "  +0.00%  "
Set-TextValue "E41" "  -1.69%  "
Set-TextValue "D42" "4.66"
Set-TextValue "E42" "  -2.95%  "
Set-TextValue "E43" "  -0.65%  "
Set-TextValue "D44" "38.09"
Set-TextValue "E44" "  -0.47%  "
Set-TextValue "E45" "  -5.96%  "
Set-TextValue "D46" "1.07"
Set-TextValue "E46" "  -8.04%  "
Set-TextValue "D47" "140.22"
Set-TextValue "E47" "  -2.23%  "
Set-TextValue "E48" "  -0.97%  "
Set-TextValue "D49" "0.516"
Set-TextValue "E49" "  -2.96%  "
Set-TextValue "E50" "  -0.78%  "
Set-TextValue "D51" "0.577"
Set-TextValue "E51" "  -0.99%  "
